$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$feed = "NIST Cybersecurity and IT news and events"

$rows = @(
    @("Webinar: #ChooseCyber - Insider Tips for Cybersecurity Career Week 2023", "https://www.nist.gov/news-events/events/2023/09/webinar-choosecyber-insider-tips-cybersecurity-career-week-2023", "2023-09-06", "Speaker: Davina Pruitt-Mentle Lead for Academic Engagement NICE Synopsis: Cybersecurity Career Week will take place on October 16-21, 2023, and you’re invited to participate! Join this webinar to learn more about #CyberCareerWeek to promote the"),
    @("NIST Drafts Major Update to Its Widely Used Cybersecurity Framework", "https://www.nist.gov/news-events/news/2023/08/nist-drafts-major-update-its-widely-used-cybersecurity-framework", "2023-08-08", "NIST has revised the framework to help benefit all sectors, not just critical infrastructure."),
    @("NIST Researchers Explore Best Practices for Talking to Kids About Online Privacy", "https://www.nist.gov/news-events/news/2023/08/nist-researchers-explore-best-practices-talking-kids-about-online-privacy", "2023-08-07", "Recent research by scientists at NIST suggests that a good strategy is to talk to your kids directly."),
    @("NIST Appoints Cherilyn Pascoe as New Director of the NCCoE", "https://www.nist.gov/news-events/news/2023/08/nist-appoints-cherilyn-pascoe-new-director-nccoe", "2023-08-01", "As NCCoE director, Cherilyn will provide strategic direction and technical leadership for the NCCoE, continue to align the NCCoE's work with the industry, government, and NIST priorities, and strengthen and grow relationships with key stakeholders"),
    @("STPPA6: Special Topics on Privacy and Public Auditability — Event 6: Community Efforts on Advanced Cryptographic Techniques", "https://www.nist.gov/news-events/events/2023/07/stppa6-special-topics-privacy-and-public-auditability-event-6-community", "2023-07-25", "In the Special Topics on Privacy and Public Auditability series, the NIST privacy-enhancing cryptography (PEC) project hosts talks on various interconnected topics related to privacy and public auditability. A main goal of the STPPA series is to")
)

$startRow = 152
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $feed
    $ws.Cells.Item($r, 2).Value = $data[0]
    $ws.Cells.Item($r, 3).Value = $data[1]
    $dc = $ws.Cells.Item($r, 4)
    $dc.NumberFormat = "@"
    $dc.Value = $data[2]
    $dc.Style = "Normal"
    $ws.Cells.Item($r, 5).Value = $data[3]
}

Write-Output "done"
